$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2210.1094
$ws.Range("J17").Value = 2243.2856
$ws.Range("L17").Value = 6729.8568
$ws.Range("N17").Value = -7065.8568
$ws.Range("H19").Value = 330.0625
$ws.Range("I19").Value = 440
$ws.Range("J19").Value = 280.0909
$ws.Range("K19").Value = 440
$ws.Range("L19").Value = 280.0909
$ws.Range("M19").Value = -265
$ws.Range("N19").Value = -630.0908999999999
$ws.Range("H62").Value = 3105.842
$ws.Range("I62").Value = 3215.4075
$ws.Range("J62").Value = 2836.9092
$ws.Range("K62").Value = 3215.4075
$ws.Range("L62").Value = 2836.9092
$ws.Range("M62").Value = -2591.4075
$ws.Range("N62").Value = -4084.9092
$ws.Range("H65").Value = 3105.842
$ws.Range("I65").Value = 3215.4075
$ws.Range("J65").Value = 2836.9092
$ws.Range("K65").Value = 16077.0375
$ws.Range("L65").Value = 14184.546
$ws.Range("M65").Value = -12957.0375
$ws.Range("N65").Value = -20424.546
$ws.Range("H97").Value = 100949.63
$ws.Range("J97").Value = 100949.63
$ws.Range("L97").Value = 302848.89
$ws.Range("N97").Value = -303840.89
$ws.Range("H137").Value = 3409.3696
$ws.Range("I137").Value = 880.9666999999999
$ws.Range("K137").Value = 2642.9001
$ws.Range("M137").Value = -92.90009999999984
$ws.Range("H141").Value = 4095.375
$ws.Range("I141").Value = 1433.5
$ws.Range("J141").Value = 7822
$ws.Range("K141").Value = 4300.5
$ws.Range("L141").Value = 23466
$ws.Range("M141").Value = 879.5
$ws.Range("N141").Value = -33826

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14280.113
$ws.Range("I32").Value = 13711.648
$ws.Range("K32").Value = 13711.648
$ws.Range("M32").Value = -13424.648
$ws.Range("H61").Value = 1209.6757
$ws.Range("I61").Value = 919.8148
$ws.Range("J61").Value = 1992.3
$ws.Range("K61").Value = 919.8148
$ws.Range("L61").Value = 1992.3
$ws.Range("M61").Value = -707.8148
$ws.Range("N61").Value = -2416.3
$ws.Range("H74").Value = 1571.4783
$ws.Range("I74").Value = 1326.5853
$ws.Range("J74").Value = 3579.6
$ws.Range("K74").Value = 1326.5853
$ws.Range("L74").Value = 3579.6
$ws.Range("M74").Value = -452.5853
$ws.Range("N74").Value = -5327.6
$ws.Range("H77").Value = 1571.4783
$ws.Range("I77").Value = 1326.5853
$ws.Range("J77").Value = 3579.6
$ws.Range("K77").Value = 6632.9265
$ws.Range("L77").Value = 17898
$ws.Range("M77").Value = -2264.9265
$ws.Range("N77").Value = -26634
$ws.Range("H132").Value = 1652.7234
$ws.Range("I132").Value = 894.67645
$ws.Range("J132").Value = 3635.3076
$ws.Range("K132").Value = 2684.02935
$ws.Range("L132").Value = 10905.9228
$ws.Range("M132").Value = -154.0293500000002
$ws.Range("N132").Value = -15965.9228
$ws.Range("H135").Value = 51447.75
$ws.Range("J135").Value = 51447.75
$ws.Range("L135").Value = 51447.75
$ws.Range("N135").Value = -61587.75
$ws.Range("H136").Value = 1209.6757
$ws.Range("I136").Value = 919.8148
$ws.Range("J136").Value = 1992.3
$ws.Range("K136").Value = 2759.4444
$ws.Range("L136").Value = 5976.9
$ws.Range("M136").Value = -209.4443999999999
$ws.Range("N136").Value = -11076.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 30260
$ws.Range("J81").Value = 30260
$ws.Range("L81").Value = 30260
$ws.Range("N81").Value = -32382
$ws.Range("H84").Value = 30260
$ws.Range("J84").Value = 30260
$ws.Range("L84").Value = 90780
$ws.Range("N84").Value = -101388
$ws.Range("H134").Value = 3281.1035
$ws.Range("I134").Value = 2646.1428
$ws.Range("J134").Value = 3483.1365
$ws.Range("K134").Value = 7938.428400000001
$ws.Range("L134").Value = 10449.4095
$ws.Range("M134").Value = -5403.428400000001
$ws.Range("N134").Value = -15519.4095
$ws.Range("H135").Value = 22803.25
$ws.Range("J135").Value = 22803.25
$ws.Range("L135").Value = 22803.25
$ws.Range("N135").Value = -32943.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4425.725
$ws.Range("I31").Value = 2568.2942
$ws.Range("K31").Value = 2568.2942
$ws.Range("M31").Value = -2273.2942
$ws.Range("H34").Value = 4425.725
$ws.Range("I34").Value = 2568.2942
$ws.Range("K34").Value = 2568.2942
$ws.Range("M34").Value = -2366.2942
$ws.Range("H58").Value = 1513.6316
$ws.Range("I58").Value = 1235.7667
$ws.Range("J58").Value = 2555.625
$ws.Range("K58").Value = 1235.7667
$ws.Range("L58").Value = 2555.625
$ws.Range("M58").Value = -1032.7667
$ws.Range("N58").Value = -2961.625
$ws.Range("H132").Value = 46170.188
$ws.Range("I132").Value = 1239.8572
$ws.Range("J132").Value = 131946.27
$ws.Range("K132").Value = 3719.5716
$ws.Range("L132").Value = 395838.8099999999
$ws.Range("M132").Value = -1189.5716
$ws.Range("N132").Value = -400898.8099999999
$ws.Range("H134").Value = 351342.3
$ws.Range("I134").Value = 963.3889
$ws.Range("J134").Value = 3504752.8
$ws.Range("K134").Value = 2890.1667
$ws.Range("L134").Value = 10514258.4
$ws.Range("M134").Value = -355.1667000000002
$ws.Range("N134").Value = -10519328.4
$ws.Range("H136").Value = 1513.6316
$ws.Range("I136").Value = 1235.7667
$ws.Range("J136").Value = 2555.625
$ws.Range("K136").Value = 3707.300099999999
$ws.Range("L136").Value = 7666.875
$ws.Range("M136").Value = -1157.300099999999
$ws.Range("N136").Value = -12766.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 29.4
$ws.Range("I14").Value = 29.4
$ws.Range("K14").Value = 88.19999999999999
$ws.Range("M14").Value = 84.80000000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2775.577
$ws.Range("I132").Value = 1560.8125
$ws.Range("J132").Value = 4719.2
$ws.Range("K132").Value = 4682.4375
$ws.Range("L132").Value = 14157.6
$ws.Range("M132").Value = -2152.4375
$ws.Range("N132").Value = -19217.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2484.8572
$ws.Range("I100").Value = 2460.6155
$ws.Range("K100").Value = 2460.6155
$ws.Range("M100").Value = -1919.6155
$ws.Range("H122").Value = 36160.277
$ws.Range("I122").Value = 51402.25
$ws.Range("J122").Value = 2289.2222
$ws.Range("K122").Value = 154206.75
$ws.Range("L122").Value = 6867.6666
$ws.Range("M122").Value = -151756.75
$ws.Range("N122").Value = -11767.6666
$ws.Range("H132").Value = 1952.0227
$ws.Range("I132").Value = 1426.0984
$ws.Range("J132").Value = 3140.2222
$ws.Range("K132").Value = 4278.2952
$ws.Range("L132").Value = 9420.6666
$ws.Range("M132").Value = -1748.2952
$ws.Range("N132").Value = -14480.6666
$ws.Range("H136").Value = 2574.6667
$ws.Range("I136").Value = 2099.5
$ws.Range("J136").Value = 2812.25
$ws.Range("K136").Value = 6298.5
$ws.Range("L136").Value = 8436.75
$ws.Range("M136").Value = -3748.5
$ws.Range("N136").Value = -13536.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1549.1063
$ws.Range("I132").Value = 1211.6285
$ws.Range("J132").Value = 2533.4167
$ws.Range("K132").Value = 3634.8855
$ws.Range("L132").Value = 7600.250100000001
$ws.Range("M132").Value = -1104.8855
$ws.Range("N132").Value = -12660.2501
$ws.Range("H136").Value = 303836
$ws.Range("I136").Value = 357763.97
$ws.Range("J136").Value = 1839.4
$ws.Range("K136").Value = 1073291.91
$ws.Range("L136").Value = 5518.200000000001
$ws.Range("M136").Value = -1070741.91
$ws.Range("N136").Value = -10618.2

Write-Output "Applied 188 cell updates across 8 sheets"